$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 98.45
$ws.Range("I15").Value = 98.45
$ws.Range("K15").Value = 295.35
$ws.Range("M15").Value = -126.35
$ws.Range("H17").Value = 11469.182
$ws.Range("J17").Value = 11469.182
$ws.Range("L17").Value = 34407.546
$ws.Range("N17").Value = -34743.546
$ws.Range("H111").Value = 2029
$ws.Range("I111").Value = 2029
$ws.Range("K111").Value = 6087
$ws.Range("M111").Value = -3020
$ws.Range("H138").Value = 2010.5638
$ws.Range("I138").Value = 1629.9048
$ws.Range("J138").Value = 2318.0193
$ws.Range("K138").Value = 4889.7144
$ws.Range("L138").Value = 6954.0579
$ws.Range("M138").Value = 250.2856000000002
$ws.Range("N138").Value = -17234.0579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3237.0217
$ws.Range("I2").Value = 3462.4524
$ws.Range("J2").Value = 870
$ws.Range("K2").Value = 3462.4524
$ws.Range("L2").Value = 870
$ws.Range("M2").Value = -3349.4524
$ws.Range("N2").Value = -1096
$ws.Range("H80").Value = 44012.8
$ws.Range("J80").Value = 49991
$ws.Range("L80").Value = 49991
$ws.Range("N80").Value = -51987
$ws.Range("H83").Value = 44012.8
$ws.Range("J83").Value = 49991
$ws.Range("L83").Value = 149973
$ws.Range("N83").Value = -159957
$ws.Range("H113").Value = 38520.125
$ws.Range("J113").Value = 38520.125
$ws.Range("L113").Value = 38520.125
$ws.Range("N113").Value = -47198.125
$ws.Range("H116").Value = 3237.0217
$ws.Range("I116").Value = 3462.4524
$ws.Range("J116").Value = 870
$ws.Range("K116").Value = 3462.4524
$ws.Range("L116").Value = 870
$ws.Range("M116").Value = -1168.4524
$ws.Range("N116").Value = -5458

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3237.0217
$ws.Range("I3").Value = 3462.4524
$ws.Range("J3").Value = 870
$ws.Range("K3").Value = 3462.4524
$ws.Range("L3").Value = 870
$ws.Range("M3").Value = -3348.4524
$ws.Range("N3").Value = -1098

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2156.9285
$ws.Range("I99").Value = 1857.8334
$ws.Range("K99").Value = 1857.8334
$ws.Range("M99").Value = -359.8334
$ws.Range("H126").Value = 2156.9285
$ws.Range("I126").Value = 1857.8334
$ws.Range("K126").Value = 5573.5002
$ws.Range("M126").Value = -3103.5002
$ws.Range("H132").Value = 31455.729
$ws.Range("I132").Value = 1368.9714
$ws.Range("J132").Value = 112458.54
$ws.Range("K132").Value = 4106.914199999999
$ws.Range("L132").Value = 337375.62
$ws.Range("M132").Value = -1576.914199999999
$ws.Range("N132").Value = -342435.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3224.26
$ws.Range("I5").Value = 5227.8096
$ws.Range("J5").Value = 1773.4138
$ws.Range("K5").Value = 15683.4288
$ws.Range("L5").Value = 5320.2414
$ws.Range("M5").Value = -15571.4288
$ws.Range("N5").Value = -5544.2414
$ws.Range("H107").Value = 10029.333
$ws.Range("I107").Value = 14817.571
$ws.Range("J107").Value = 7635.2144
$ws.Range("K107").Value = 44452.713
$ws.Range("L107").Value = 22905.6432
$ws.Range("M107").Value = -42532.713
$ws.Range("N107").Value = -26745.6432
$ws.Range("H113").Value = 7357.3335
$ws.Range("I113").Value = 7847.143
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 23541.429
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -21371.429
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 2597.8076
$ws.Range("I122").Value = 616.2917
$ws.Range("J122").Value = 26376
$ws.Range("K122").Value = 5546.6253
$ws.Range("L122").Value = 237384
$ws.Range("M122").Value = -3096.6253
$ws.Range("N122").Value = -242284
$ws.Range("H132").Value = 1156.8125
$ws.Range("I132").Value = 880.4
$ws.Range("J132").Value = 1617.5
$ws.Range("K132").Value = 7923.599999999999
$ws.Range("L132").Value = 14557.5
$ws.Range("M132").Value = -5393.599999999999
$ws.Range("N132").Value = -19617.5
$ws.Range("H135").Value = 3224.26
$ws.Range("I135").Value = 5227.8096
$ws.Range("J135").Value = 1773.4138
$ws.Range("K135").Value = 47050.2864
$ws.Range("L135").Value = 15960.7242
$ws.Range("M135").Value = -44515.2864
$ws.Range("N135").Value = -21030.7242
$ws.Range("H139").Value = 5362.758
$ws.Range("I139").Value = 6534.591
$ws.Range("J139").Value = 3019.0908
$ws.Range("K139").Value = 19603.773
$ws.Range("L139").Value = 9057.2724
$ws.Range("M139").Value = -14463.773
$ws.Range("N139").Value = -19337.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 242531.72
$ws.Range("I80").Value = 389499.22
$ws.Range("J80").Value = 3709.5
$ws.Range("K80").Value = 389499.22
$ws.Range("L80").Value = 3709.5
$ws.Range("M80").Value = -388501.22
$ws.Range("N80").Value = -5705.5
$ws.Range("H83").Value = 242531.72
$ws.Range("I83").Value = 389499.22
$ws.Range("J83").Value = 3709.5
$ws.Range("K83").Value = 1947496.1
$ws.Range("L83").Value = 18547.5
$ws.Range("M83").Value = -1942504.1
$ws.Range("N83").Value = -28531.5
$ws.Range("H132").Value = 3377.8647
$ws.Range("I132").Value = 3225
$ws.Range("K132").Value = 9675
$ws.Range("M132").Value = -7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1350.2727
$ws.Range("I16").Value = 1559.2222
$ws.Range("J16").Value = 410
$ws.Range("K16").Value = 1559.2222
$ws.Range("L16").Value = 410
$ws.Range("M16").Value = -1389.2222
$ws.Range("N16").Value = -750
$ws.Range("H22").Value = 1103.55
$ws.Range("I22").Value = 1136.1538
$ws.Range("J22").Value = 1043
$ws.Range("K22").Value = 1136.1538
$ws.Range("L22").Value = 1043
$ws.Range("M22").Value = -841.1538
$ws.Range("N22").Value = -1633
$ws.Range("H27").Value = 1103.55
$ws.Range("I27").Value = 1136.1538
$ws.Range("J27").Value = 1043
$ws.Range("K27").Value = 1136.1538
$ws.Range("L27").Value = 1043
$ws.Range("M27").Value = -1029.1538
$ws.Range("N27").Value = -1257
$ws.Range("H93").Value = 1960
$ws.Range("J93").Value = 3990
$ws.Range("L93").Value = 3990
$ws.Range("N93").Value = -6486
$ws.Range("H100").Value = 2662.375
$ws.Range("I100").Value = 2450.8333
$ws.Range("K100").Value = 2450.8333
$ws.Range("M100").Value = -1909.8333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 50002104
$ws.Range("J107").Value = 50002104
$ws.Range("L107").Value = 150006312
$ws.Range("N107").Value = -150010152
$ws.Range("H136").Value = 16047.776
$ws.Range("I136").Value = 29021.691
$ws.Range("J136").Value = 2372.5676
$ws.Range("K136").Value = 87065.073
$ws.Range("L136").Value = 7117.702799999999
$ws.Range("M136").Value = -84515.073
$ws.Range("N136").Value = -12217.7028
